$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: becomes the "subject" row ---
$ws.Range("B2").Value = "subject"
$ws.Range("C2").Value = "Activate authentication via Gmail"
$ws.Range("E2").Value = 44588.0

# --- Row 3: becomes the "letterName" row ---
$ws.Range("B3").Value = "letterName"
$ws.Range("C3").Value = "Activate authentication via Gmail"
$ws.Range("E3").Value = 44588.0

# Copy the "Updated By" text value (keeps it as text "960024377", same as rows 4-8)
# into D2/D3, preserving their existing cell format.
$ws.Range("D4").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("D4").Copy()
$ws.Range("D3").PasteSpecial(-4163)

# Apply the same date number formatting used by the other "Last Updated" cells
# (E4:E8) to the newly-populated E2/E3 date values.
$ws.Range("E4").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E3").PasteSpecial(-4122)

$excel.CutCopyMode = $false
